# Automatische sync: 2025-06-17 20:41:37
# Appends two new mail-log entries (rows 21 and 22) to the "Logs" sheet and
# refreshes the "Informatieaanvraag" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")

# --- Row 21 -----------------------------------------------------------
$logs.Range("A21").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B21").Value = "mailmind.test@zohomail.eu"
$logs.Range("C21").Value = "Hallo, ik wil graag weten wanneer jullie bereikbaar zijn.`nSent using {0}"
$logs.Range("D21").Value = "Informatieaanvraag"
$logs.Range("E21").Value = "Beste afzender,`nBedankt voor je e-mail. Onze kantooruren zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Buiten deze tijden proberen we zo snel mogelijk te reageren. Mocht je verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[E-mailassistent]"
$logs.Range("F21").Value = "2025-06-17 20:40:26"
$logs.Range("G21").Value = "Ja"

# --- Row 22 -----------------------------------------------------------
$logs.Range("A22").Value = "Re: Wat zijn jullie openingstijden?"
$logs.Range("B22").Value = "mailmind.test@zohomail.eu"
$logs.Range("C22").Value = "Beste afzender,`nBedankt voor je e-mail. Onze kantooruren zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Buiten deze tijden proberen we zo snel mogelijk te reageren. Mocht je verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[E-mailassistent]"
$logs.Range("D22").Value = "Informatieaanvraag"
$logs.Range("E22").Value = "Beste afzender,`nDank voor je bericht. Onze kantooruren zijn van maandag t/m vrijdag van 9:00-17:00 uur. Buiten deze tijden streven we naar spoedige reactie. Voor verdere vragen, neem gerust contact met ons op.`nMet vriendelijke groet,`n[E-mailassistent]"
$logs.Range("F22").Value = "2025-06-17 20:40:38"
$logs.Range("G22").Value = "Ja"

# Re-fit the row heights so the new rows keep the sheet's default height
# instead of Excel's "wrap multi-line text" auto-expansion.
$logs.Range("A21:G22").Rows.AutoFit()

# --- Extend conditional formatting to cover the new rows ---------------
$fcsCategory = $logs.Range("D2:D20").FormatConditions
for ($i = 1; $i -le $fcsCategory.Count; $i++) {
    $fcsCategory.Item($i).ModifyAppliesToRange($logs.Range("D2:D22"))
}

$fcsAnswered = $logs.Range("G2:G20").FormatConditions
for ($i = 1; $i -le $fcsAnswered.Count; $i++) {
    $fcsAnswered.Item($i).ModifyAppliesToRange($logs.Range("G2:G22"))
}

# --- Dashboard tally ----------------------------------------------------
# Two more "Informatieaanvraag" mails came in, so the count goes from 7 to 9.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 9
